$wb = $excel.ActiveWorkbook

# Data to append to each sheet: two new rows (dates 45967, 45968) with remn_amt values.
# Sheet order matches workbook sheet order: 카카오, NAVER, 농심, 엔씨소프트
$newData = @{
    1 = @(@(45967, 832079), @(45968, 0))
    2 = @(@(45967, 1317194), @(45968, 0))
    3 = @(@(45967, 125092), @(45968, 0))
    4 = @(@(45967, 158092), @(45968, 0))
}

foreach ($sheetIndex in 1..4) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $rows = $newData[$sheetIndex]
    $startRow = 104
    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $startRow + $i
        $dateVal = $rows[$i][0]
        $amtVal = $rows[$i][1]

        $cellA = $ws.Cells.Item($r, 1)
        $cellA.Value = $dateVal
        $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"

        $cellB = $ws.Cells.Item($r, 2)
        $cellB.Value = $amtVal
    }
}
